{"js": "// Update the worksheet date and the 25 division problems/answers.\n//\n// The document body (in document order) is: a title paragraph with the\n// date, followed by a 5-column table whose cells hold one division\n// expression each (interspersed with several fully empty rows). We walk\n// every paragraph in the body (this also reaches into table cells) and\n// replace the text of each non-empty paragraph with its new value, in\n// document order, leaving every empty paragraph (and all formatting)\n// untouched.\n\nconst oldTexts = [\n  \"2024-05-19 Sunday\",\n  \"45\u00f74=11, 1\",\n  \"17\u00f72=8, 1\",\n  \"38\u00f75=7, 3\",\n  \"54\u00f75=10, 4\",\n  \"70\u00f72=35, 0\",\n  \"11\u00f74=2, 3\",\n  \"24\u00f75=4, 4\",\n  \"73\u00f74=18, 1\",\n  \"94\u00f76=15, 4\",\n  \"28\u00f74=7, 0\",\n  \"47\u00f78=5, 7\",\n  \"37\u00f75=7, 2\",\n  \"25\u00f72=12, 1\",\n  \"39\u00f79=4, 3\",\n  \"68\u00f73=22, 2\",\n  \"43\u00f74=10, 3\",\n  \"64\u00f73=21, 1\",\n  \"35\u00f76=5, 5\",\n  \"12\u00f74=3, 0\",\n  \"91\u00f75=18, 1\",\n  \"80\u00f77=11, 3\",\n  \"94\u00f76=15, 4\",\n  \"65\u00f79=7, 2\",\n  \"31\u00f72=15, 1\",\n  \"82\u00f73=27, 1\",\n];\n\nconst newTexts = [\n  \"2024-05-20 Monday\",\n  \"96\u00f79=10, 6\",\n  \"17\u00f75=3, 2\",\n  \"89\u00f77=12, 5\",\n  \"41\u00f77=5, 6\",\n  \"53\u00f78=6, 5\",\n  \"79\u00f76=13, 1\",\n  \"37\u00f75=7, 2\",\n  \"51\u00f73=17, 0\",\n  \"80\u00f79=8, 8\",\n  \"48\u00f78=6, 0\",\n  \"70\u00f74=17, 2\",\n  \"14\u00f73=4, 2\",\n  \"68\u00f77=9, 5\",\n  \"11\u00f76=1, 5\",\n  \"20\u00f79=2, 2\",\n  \"26\u00f75=5, 1\",\n  \"24\u00f78=3, 0\",\n  \"33\u00f76=5, 3\",\n  \"15\u00f75=3, 0\",\n  \"80\u00f77=11, 3\",\n  \"35\u00f79=3, 8\",\n  \"14\u00f79=1, 5\",\n  \"50\u00f74=12, 2\",\n  \"70\u00f79=7, 7\",\n  \"82\u00f77=11, 5\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet matchIndex = 0;\nfor (let i = 0; i < paragraphs.items.length && matchIndex < newTexts.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  if (text === \"\") continue; // skip the blank spacer rows in the table\n\n  if (text !== oldTexts[matchIndex]) {\n    throw new Error(\n      `Unexpected paragraph text at position ${matchIndex}: got ${JSON.stringify(\n        text\n      )}, expected ${JSON.stringify(oldTexts[matchIndex])}`\n    );\n  }\n\n  para.insertText(newTexts[matchIndex], \"Replace\");\n  matchIndex++;\n}\n\nawait context.sync();\n\nif (matchIndex !== newTexts.length) {\n  throw new Error(`Only replaced ${matchIndex} of ${newTexts.length} expected paragraphs`);\n}\n", "ps1": "# Update the worksheet date and the 25 division problems/answers.\n#\n# The document (in paragraph order) is: a title paragraph with the date,\n# followed by a 5-column table whose cells hold one division expression\n# each (interspersed with several fully empty rows). $d.Paragraphs walks\n# every paragraph in document order, including the ones inside table\n# cells, so we can address each target paragraph by its 1-based COM\n# index and overwrite just its Range.Text -- this keeps every run's\n# formatting (fonts/size) untouched and leaves every blank cell alone.\n\n$d = $word.ActiveDocument\n\n$targets = @(\n    @{ Index = 1;   Old = \"2024-05-19 Sunday\"; New = \"2024-05-20 Monday\" }\n    @{ Index = 2;   Old = \"45\u00f74=11, 1\";         New = \"96\u00f79=10, 6\" }\n    @{ Index = 3;   Old = \"17\u00f72=8, 1\";          New = \"17\u00f75=3, 2\" }\n    @{ Index = 4;   Old = \"38\u00f75=7, 3\";          New = \"89\u00f77=12, 5\" }\n    @{ Index = 5;   Old = \"54\u00f75=10, 4\";         New = \"41\u00f77=5, 6\" }\n    @{ Index = 6;   Old = \"70\u00f72=35, 0\";         New = \"53\u00f78=6, 5\" }\n    @{ Index = 26;  Old = \"11\u00f74=2, 3\";          New = \"79\u00f76=13, 1\" }\n    @{ Index = 27;  Old = \"24\u00f75=4, 4\";          New = \"37\u00f75=7, 2\" }\n    @{ Index = 28;  Old = \"73\u00f74=18, 1\";         New = \"51\u00f73=17, 0\" }\n    @{ Index = 29;  Old = \"94\u00f76=15, 4\";         New = \"80\u00f79=8, 8\" }\n    @{ Index = 30;  Old = \"28\u00f74=7, 0\";          New = \"48\u00f78=6, 0\" }\n    @{ Index = 50;  Old = \"47\u00f78=5, 7\";          New = \"70\u00f74=17, 2\" }\n    @{ Index = 51;  Old = \"37\u00f75=7, 2\";          New = \"14\u00f73=4, 2\" }\n    @{ Index = 52;  Old = \"25\u00f72=12, 1\";         New = \"68\u00f77=9, 5\" }\n    @{ Index = 53;  Old = \"39\u00f79=4, 3\";          New = \"11\u00f76=1, 5\" }\n    @{ Index = 54;  Old = \"68\u00f73=22, 2\";         New = \"20\u00f79=2, 2\" }\n    @{ Index = 74;  Old = \"43\u00f74=10, 3\";         New = \"26\u00f75=5, 1\" }\n    @{ Index = 75;  Old = \"64\u00f73=21, 1\";         New = \"24\u00f78=3, 0\" }\n    @{ Index = 76;  Old = \"35\u00f76=5, 5\";          New = \"33\u00f76=5, 3\" }\n    @{ Index = 77;  Old = \"12\u00f74=3, 0\";          New = \"15\u00f75=3, 0\" }\n    @{ Index = 78;  Old = \"91\u00f75=18, 1\";         New = \"80\u00f77=11, 3\" }\n    @{ Index = 98;  Old = \"80\u00f77=11, 3\";         New = \"35\u00f79=3, 8\" }\n    @{ Index = 99;  Old = \"94\u00f76=15, 4\";         New = \"14\u00f79=1, 5\" }\n    @{ Index = 100; Old = \"65\u00f79=7, 2\";          New = \"50\u00f74=12, 2\" }\n    @{ Index = 101; Old = \"31\u00f72=15, 1\";         New = \"70\u00f79=7, 7\" }\n    @{ Index = 102; Old = \"82\u00f73=27, 1\";         New = \"82\u00f77=11, 5\" }\n)\n\nforeach ($t in $targets) {\n    $p = $d.Paragraphs($t.Index)\n    $r = $p.Range\n    $current = $r.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $t.Old) {\n        throw \"Paragraph $($t.Index) text mismatch: got '$current', expected '$($t.Old)'\"\n    }\n    $r.Text = $t.New\n}\n"}
